$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows -----------------------------------------------
# New row 2: "Aguascalientes" record (Rf7971) goes above the current row 2.
$ws.Rows.Item(2).Insert()
# New row (now 4, after the shift above): "Gastronomia BJ" record (R56a72)
# is inserted right after the first "Baja California" / R23387 row (row 3).
$ws.Rows.Item(4).Insert()

# --- Helper to write a full data row ------------------------------------
function Set-DataRow {
    param(
        [int]$Row,
        [string]$RecordId,
        [string]$AuthorId,
        [string]$StateName,
        [string]$RecordType,
        [bool]$IsPublic,
        [string]$Title,
        [string]$Description,
        [string]$ImageUrl,
        [bool]$Deleted
    )
    $ws.Cells.Item($Row, 1).Value = $RecordId
    $ws.Cells.Item($Row, 2).Value = $AuthorId
    $ws.Cells.Item($Row, 3).Value = $StateName
    $ws.Cells.Item($Row, 4).Value = $RecordType
    $ws.Cells.Item($Row, 5).Value = $IsPublic
    $ws.Cells.Item($Row, 6).Value = $Title
    $ws.Cells.Item($Row, 7).Value = $Description
    $ws.Cells.Item($Row, 8).Value = $ImageUrl
    $ws.Cells.Item($Row, 9).Value = $Deleted
}

# Row 2 (new): Aguascalientes
Set-DataRow 2 'Rf7971' 'Udd529' 'Aguascalientes' 'Gastronomía' $true 'Aguascalientes' 'Aguascalientes' 'src/main/java/proyecto/resources/agus/489496123_1872564096894029_4907360869533054311_n.jpg' $false

# Row 3 (was row 2, unchanged content): R23387 / Baja California
Set-DataRow 3 'R23387' 'Udd529' 'Baja California' 'Tradición' $true 'asdsada' 'adsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasdadsadasdasd' 'src/main/java/proyecto/resources/baca/download.jpg' $false

# Row 4 (new): Gastronomia BJ
Set-DataRow 4 'R56a72' 'Udd529' 'Baja California' 'Gastronomía' $true 'Gastronomia BJ' 'Gastronomia BJ' 'src/main/java/proyecto/resources/baca/492522295_1240113004212486_7120062313825821510_n.jpg' $false

# Row 5 (was row 3, unchanged content): Rc35d7 / Baja California / Lugar
Set-DataRow 5 'Rc35d7' 'Udd529' 'Baja California' 'Lugar' $true 'adasd' 'dsadasda' 'src/main/java/proyecto/resources/baca/imgnotfound.png' $false

# Row 6 (was row 4, unchanged content): Re2704 / Campeche / Tradicion
Set-DataRow 6 'Re2704' 'Udd529' 'Campeche' 'Tradición' $true 'rrruuu' 'adsasdasdasd' 'src/main/java/proyecto/resources/camp/460517427_1926744817802502_6764676756989845370_n.jpg' $false

# Row 7 (was row 5, unchanged content): R4e977 / Campeche / Lugar
Set-DataRow 7 'R4e977' 'Uf197a' 'Campeche' 'Lugar' $true 'Prueba de con backup #2' 'Desproporcion' 'src/main/java/proyecto/resources/camp/483065467_10237970398099334_6612370661461858019_n.jpg' $false

# Row 8 (was row 6, unchanged content): R95be5 / Campeche / Regionalismo o localismo
Set-DataRow 8 'R95be5' 'Uf197a' 'Campeche' 'Regionalismo o localismo' $false 'Prueba sin eliminacion' 'dasjhdiahdkjsa' 'src/main/java/proyecto/resources/camp/555472471_776962455177966_3128049131977084705_n.jpg' $false

# Row 9 (was row 7, unchanged content): R52ec2 / Puebla / Gastronomia
Set-DataRow 9 'R52ec2' 'Uf197a' 'Puebla' 'Gastronomía' $false 'Prueba de eliminacion ' 'adasdsad' 'src/main/java/proyecto/resources/pueb/474921981_927529119567042_1843877378970278987_n.jpg' $true

